$d = $word.ActiveDocument

# Insert a new bold Arial 20pt paragraph "cincong 1" before the current
# first paragraph ("sf sdsf sdf dsf"). InsertParagraphBefore on the first
# paragraph's range creates a new empty paragraph ahead of it, inheriting
# the run formatting of what follows, so we just need to set its text.
$firstRange = $d.Paragraphs.First.Range
$firstRange.InsertParagraphBefore()
$d.Paragraphs.First.Range.Text = "cincong 1"

# Append a new plain paragraph after the last (Lorem Ipsum) paragraph.
$lastRange = $d.Paragraphs.Last.Range
$lastRange.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "sf sdsf sdf dsfsf sdsf sdf dsfsf sdsf sdf dsf"
